# Updates cryptos list prices/volume values per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, [string]$Address, [string]$NewValue)
    $cell = $Worksheet.Range($Address)
    # Leading apostrophe forces Excel to store the value as literal text
    # (numeric-looking strings like '1.001' or '0.000007311' would otherwise
    # be auto-converted to numbers).
    $cell.Value = "'" + $NewValue
    # Re-apply the Normal style so no stray number-format/quote-prefix style
    # sticks to the cell (keeps the cell style identical to before the edit).
    $cell.Style = "Normal"
}

Set-TextValue $ws 'D2' '30.521.32'
Set-TextValue $ws 'E2' '  +0.98%  '
Set-TextValue $ws 'D3' '1.852.68'
Set-TextValue $ws 'E3' '  +0.19%  '
Set-TextValue $ws 'D4' '1.001'
Set-TextValue $ws 'D5' '233.56'
Set-TextValue $ws 'E5' '  +0.23%  '
Set-TextValue $ws 'D6' '1.001'
Set-TextValue $ws 'E6' '  +0.04%  '
Set-TextValue $ws 'D7' '0.4741'
Set-TextValue $ws 'E7' '  +0.44%  '
Set-TextValue $ws 'D8' '0.2742'
Set-TextValue $ws 'E8' '  +0.88%  '
Set-TextValue $ws 'D9' '0.06305'
Set-TextValue $ws 'E9' '  -1.09%  '
Set-TextValue $ws 'D10' '17.67'
Set-TextValue $ws 'E10' '  +8.97%  '
Set-TextValue $ws 'D11' '1.853.16'
Set-TextValue $ws 'E11' '  +0.19%  '
Set-TextValue $ws 'D12' '0.07445'
Set-TextValue $ws 'E12' '  +0.30%  '
Set-TextValue $ws 'D13' '4.976'
Set-TextValue $ws 'E13' '  +0.88%  '
Set-TextValue $ws 'D14' '84.53'
Set-TextValue $ws 'E14' '  -0.62%  '
Set-TextValue $ws 'D15' '0.6256'
Set-TextValue $ws 'E15' '  -0.28%  '
Set-TextValue $ws 'D16' '30.482.82'
Set-TextValue $ws 'E16' '  +1.03%  '
Set-TextValue $ws 'D17' '245.92'
Set-TextValue $ws 'E17' '  +8.44%  '
Set-TextValue $ws 'D18' '1.000'
Set-TextValue $ws 'E18' '  +0.09%  '
Set-TextValue $ws 'D19' '12.68'
Set-TextValue $ws 'E19' '  +0.61%  '
Set-TextValue $ws 'D20' '0.000007311'
Set-TextValue $ws 'E20' '  -0.32%  '
Set-TextValue $ws 'D21' '1.000'
Set-TextValue $ws 'E21' '  -0.10%  '
Set-TextValue $ws 'D22' '4.937'
Set-TextValue $ws 'E22' '  +0.29%  '
Set-TextValue $ws 'D23' '5.915'
Set-TextValue $ws 'E23' '  -0.15%  '
Set-TextValue $ws 'D24' '9.116'
Set-TextValue $ws 'E24' '  -1.08%  '
Set-TextValue $ws 'D25' '162.69'
Set-TextValue $ws 'E25' '  -2.65%  '
Set-TextValue $ws 'D26' '17.97'
Set-TextValue $ws 'E26' '  +0.03%  '
Set-TextValue $ws 'D27' '1.873'
Set-TextValue $ws 'E27' '  +0.00%  '
Set-TextValue $ws 'E28' '  +0.11%  '
Set-TextValue $ws 'D29' '1.353'
Set-TextValue $ws 'E29' '  -1.93%  '
Set-TextValue $ws 'D30' '4.006'
Set-TextValue $ws 'E30' '  -2.56%  '
Set-TextValue $ws 'D31' '3.826'
Set-TextValue $ws 'E31' '  -1.53%  '
Set-TextValue $ws 'D32' '0.04838'
Set-TextValue $ws 'E32' '  -0.91%  '
Set-TextValue $ws 'D33' '1.133'
Set-TextValue $ws 'E33' '  -1.89%  '
Set-TextValue $ws 'D34' '0.7016'
Set-TextValue $ws 'E34' '  -1.18%  '
Set-TextValue $ws 'D35' '2.702'
Set-TextValue $ws 'E35' '  +0.09%  '
Set-TextValue $ws 'D36' '0.01895'
Set-TextValue $ws 'E36' '  +2.89%  '
Set-TextValue $ws 'D37' '2.688'
Set-TextValue $ws 'E37' '  +2.22%  '
Set-TextValue $ws 'D38' '1.997'
Set-TextValue $ws 'E38' '  +2.35%  '
Set-TextValue $ws 'D39' '0.8750'
Set-TextValue $ws 'E39' '  -2.97%  '
Set-TextValue $ws 'D40' '106.73'
Set-TextValue $ws 'E40' '  +1.78%  '
Set-TextValue $ws 'E41' '  +0.27%  '
Set-TextValue $ws 'D42' '5.549'
Set-TextValue $ws 'E42' '  +0.08%  '
Set-TextValue $ws 'D43' '0.4051'
Set-TextValue $ws 'E43' '  -0.67%  '
Set-TextValue $ws 'D44' '7.181'
Set-TextValue $ws 'E44' '  +1.65%  '
Set-TextValue $ws 'D45' '62.79'
Set-TextValue $ws 'E45' '  +4.33%  '
Set-TextValue $ws 'D46' '0.1203'
Set-TextValue $ws 'E46' '  +0.88%  '
Set-TextValue $ws 'D47' '33.51'
Set-TextValue $ws 'E47' '  +1.12%  '
Set-TextValue $ws 'D48' '8.564'
Set-TextValue $ws 'E48' '  -0.42%  '
Set-TextValue $ws 'D49' '0.05532'
Set-TextValue $ws 'E49' '  -0.54%  '
Set-TextValue $ws 'D50' '1.349'
Set-TextValue $ws 'E50' '  -2.52%  '
Set-TextValue $ws 'D51' '0.3688'
Set-TextValue $ws 'E51' '  +0.11%  '
